# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column E ("municipio-nombre") metadata rows are updated to use the
# curated sdmx-dimension:refArea dimension instead of the old
# iaest-measure:municipio-nombre measure, mirroring columns G and J
# (provincia-nombre / comarca-nombre) which already use this pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("E3").Value = "dim"
$ws.Range("E4").Value = "URI-Municipio"
